# Weekly crypto price/volume refresh (GitHub Actions data pull)
# Also: row 48 gains a new entry (BabyDogeCoin), shifting TheSandbox,
# RenderToken and Algorand down by one row; Cronos drops off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.365.92"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.80"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.16"
$ws.Range("E5").Value = "  -1.21%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6356"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8: Dogecoin
$ws.Range("E8").Value = "  -0.10%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  -0.79%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.61"
$ws.Range("E10").Value = "  +1.74%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07731"
$ws.Range("E11").Value = "  +0.58%  "

# Row 12: WrappedEther
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.52"
$ws.Range("E12").Value = "  -1.22%  "

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6832"
$ws.Range("E14").Value = "  -0.49%  "

# Row 15: Litecoin
$ws.Range("E15").Value = "  -1.00%  "

# Row 16: ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009932"
$ws.Range("E16").Value = "  +2.08%  "

# Row 17: Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.164"
$ws.Range("E17").Value = "  -0.98%  "

# Row 18: WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.384.22"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.84"
$ws.Range("E19").Value = "  -2.83%  "

# Row 20: Avalanche
$ws.Range("E20").Value = "  -0.66%  "

# Row 21: Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22: Chainlink
$ws.Range("E22").Value = "  -0.92%  "

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.01"
$ws.Range("E24").Value = "  +0.71%  "

# Row 25: Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1402"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.368"
$ws.Range("E26").Value = "  -1.03%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.65"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28: PancakeSwap
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.464"
$ws.Range("E28").Value = "  -1.60%  "

# Row 29: Hedera
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05704"
$ws.Range("E29").Value = "  -2.71%  "

# Row 30: Toncoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.247"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31: Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  +0.02%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33: LidoDAOToken
$ws.Range("E33").Value = "  -2.96%  "

# Row 34: ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.156"
$ws.Range("E34").Value = "  -1.23%  "

# Row 35: ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7166"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36: HuobiToken
$ws.Range("E36").Value = "  -0.19%  "

# Row 37: Maker
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.254.40"
$ws.Range("E37").Value = "  +1.41%  "

# Row 38: MXToken
$ws.Range("E38").Value = "  -0.33%  "

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01805"
$ws.Range("E39").Value = "  +1.76%  "

# Row 40: TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9076"
$ws.Range("E40").Value = "  -0.76%  "

# Row 41: FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.188"
$ws.Range("E41").Value = "  +1.17%  "

# Row 42: PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43: RocketPoolETH
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.004.24"
$ws.Range("E43").Value = "  -1.28%  "

# Row 44: Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.75"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45: Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.45"
$ws.Range("E45").Value = "  -1.55%  "

# Row 46: EnergySwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.190"
$ws.Range("E46").Value = "  +0.47%  "

# Row 47: Aptos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.049"
$ws.Range("E47").Value = "  -4.08%  "

# Row 48: BabyDogeCoin
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000117"
$ws.Range("E48").Value = "  -0.59%  "

# Row 49: TheSandbox
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4025"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50: RenderToken
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.704"
$ws.Range("E50").Value = "  +0.49%  "

# Row 51: Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1128"
$ws.Range("E51").Value = "  +0.81%  "
